# Commit #13 - Halved FM input range
# - Row 21 (47K resistors): remove R2,R5,R48 from designator list, keep R19,R22; Qty 5 -> 2
# - Row 23 (100K resistors): add R2,R5,R48 to designator list; Qty 14 -> 17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "R19,R22"
$ws.Range("E21").Value = 2

$ws.Range("A23").Value = "R1,R2,R3,R4,R5,R6,R7,R9,R10,R13,R16,R17,R20,R37,R42,R48,R49"
$ws.Range("E23").Value = 17
